# Chapter on JDBC updated.
# Slide 12 ("Exercise: Library"), Content Placeholder 2, first paragraph:
#   "Create a library application where a user using a console interface is able to:"
# becomes
#   "Create a library application where using a console interface a user is able to:"
# i.e. the phrase "a user " is moved from right after "where " to right before "is able to:".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# 1) Pull the phrase "a user " (chars 36-42 of the original text) forward so it sits
#    immediately in front of "is " (chars 69-71), producing a new run "a user is ".
$isPart = $tr.Characters(69, 3)
$isPart.Text = "a user " + $isPart.Text

# 2) Remove the now-duplicated "a user " left behind in its original spot.
$dup = $tr.Characters(36, 7)
$dup.Text = ""

# 3) Re-write "using " in place so it becomes its own run, separate from
#    "a console interface " that follows it.
$using = $tr.Characters(36, 6)
$using.Text = $using.Text

# 4) Merge "a library application " with the following "where " into a single run.
$merged = $tr.Characters(8, 28)
$merged.Text = $merged.Text
